$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 96.320746
$ws.Range("H2").Value = 288.962238
$ws.Range("I2").Value = 0.3809824610908788
$ws.Range("J2").Value = 0.3809824610908788
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 28.21795633333333
$ws.Range("N2").Value = 84.653869
$ws.Range("O2").Value = 0.2124568395711989
$ws.Range("P2").Value = 0.212456839571199
$ws.Range("Q2").Value = 2717.974604622091
$ws.Range("R2").Value = 24461.77144159882
$ws.Range("S2").Value = 0.08094232961542537
$ws.Range("T2").Value = 0.08094232961542538

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 96.320746
$ws.Range("H3").Value = 288.962238
$ws.Range("I3").Value = 0.3809824610908788
$ws.Range("J3").Value = 0.3809824610908788
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 63.95730733333333
$ws.Range("N3").Value = 191.871922
$ws.Range("O3").Value = 0.4815432848151522
$ws.Range("P3").Value = 0.4815432848151524
$ws.Range("Q3").Value = 6160.415554497937
$ws.Range("R3").Value = 55443.73999048144
$ws.Range("S3").Value = 0.1834595457706627
$ws.Range("T3").Value = 0.1834595457706628

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 96.320746
$ws.Range("H4").Value = 288.962238
$ws.Range("I4").Value = 0.3809824610908788
$ws.Range("J4").Value = 0.3809824610908788
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 27.47719633333334
$ws.Range("N4").Value = 82.431589
$ws.Range("O4").Value = 0.2068795565595709
$ws.Range("P4").Value = 0.2068795565595709
$ws.Range("Q4").Value = 2646.624048815132
$ws.Range("R4").Value = 23819.61643933618
$ws.Range("S4").Value = 0.07881748260745498
$ws.Range("T4").Value = 0.078817482607455

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 96.320746
$ws.Range("H5").Value = 288.962238
$ws.Range("I5").Value = 0.3809824610908788
$ws.Range("J5").Value = 0.3809824610908788
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.164899
$ws.Range("N5").Value = 39.494697
$ws.Range("O5").Value = 0.09912031905407785
$ws.Range("P5").Value = 0.0991203190540779
$ws.Range("Q5").Value = 1268.052892694654
$ws.Range("R5").Value = 11412.47603425189
$ws.Range("S5").Value = 0.03776310309733571
$ws.Range("T5").Value = 0.03776310309733572

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 18.46467533333333
$ws.Range("H6").Value = 55.394026
$ws.Range("I6").Value = 0.07303429161291354
$ws.Range("J6").Value = 0.07303429161291354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 28.21795633333333
$ws.Range("N6").Value = 84.653869
$ws.Range("O6").Value = 0.2124568395711989
$ws.Range("P6").Value = 0.212456839571199
$ws.Range("Q6").Value = 521.035402265177
$ws.Range("R6").Value = 4689.318620386593
$ws.Range("S6").Value = 0.01551663477640093
$ws.Range("T6").Value = 0.01551663477640093

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 18.46467533333333
$ws.Range("H7").Value = 55.394026
$ws.Range("I7").Value = 0.07303429161291354
$ws.Range("J7").Value = 0.07303429161291354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 63.95730733333333
$ws.Range("N7").Value = 191.871922
$ws.Range("O7").Value = 0.4815432848151522
$ws.Range("P7").Value = 0.4815432848151524
$ws.Range("Q7").Value = 1180.950915104219
$ws.Range("R7").Value = 10628.55823593797
$ws.Range("S7").Value = 0.03516917268743011
$ws.Range("T7").Value = 0.03516917268743012

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 18.46467533333333
$ws.Range("H8").Value = 55.394026
$ws.Range("I8").Value = 0.07303429161291354
$ws.Range("J8").Value = 0.07303429161291354
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 27.47719633333334
$ws.Range("N8").Value = 82.431589
$ws.Range("O8").Value = 0.2068795565595709
$ws.Range("P8").Value = 0.2068795565595709
$ws.Range("Q8").Value = 507.3575093652571
$ws.Range("R8").Value = 4566.217584287314
$ws.Range("S8").Value = 0.01510930186252194
$ws.Range("T8").Value = 0.01510930186252194

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 18.46467533333333
$ws.Range("H9").Value = 55.394026
$ws.Range("I9").Value = 0.07303429161291354
$ws.Range("J9").Value = 0.07303429161291354
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.164899
$ws.Range("N9").Value = 39.494697
$ws.Range("O9").Value = 0.09912031905407785
$ws.Range("P9").Value = 0.0991203190540779
$ws.Range("Q9").Value = 243.0855858311247
$ws.Range("R9").Value = 2187.770272480122
$ws.Range("S9").Value = 0.007239182286560552
$ws.Range("T9").Value = 0.007239182286560555

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 122.909391
$ws.Range("H10").Value = 368.728173
$ws.Range("I10").Value = 0.4861499128584522
$ws.Range("J10").Value = 0.4861499128584522
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 28.21795633333333
$ws.Range("N10").Value = 84.653869
$ws.Range("O10").Value = 0.2124568395711989
$ws.Range("P10").Value = 0.212456839571199
$ws.Range("Q10").Value = 3468.251828194593
$ws.Range("R10").Value = 31214.26645375133
$ws.Range("S10").Value = 0.1032858740437205
$ws.Range("T10").Value = 0.1032858740437205

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 122.909391
$ws.Range("H11").Value = 368.728173
$ws.Range("I11").Value = 0.4861499128584522
$ws.Range("J11").Value = 0.4861499128584522
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 63.95730733333333
$ws.Range("N11").Value = 191.871922
$ws.Range("O11").Value = 0.4815432848151522
$ws.Range("P11").Value = 0.4815432848151524
$ws.Range("Q11").Value = 7860.953694339833
$ws.Range("R11").Value = 70748.58324905849
$ws.Range("S11").Value = 0.2341022259504591
$ws.Range("T11").Value = 0.2341022259504591

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 122.909391
$ws.Range("H12").Value = 368.728173
$ws.Range("I12").Value = 0.4861499128584522
$ws.Range("J12").Value = 0.4861499128584522
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 27.47719633333334
$ws.Range("N12").Value = 82.431589
$ws.Range("O12").Value = 0.2068795565595709
$ws.Range("P12").Value = 0.2068795565595709
$ws.Range("Q12").Value = 3377.205467717433
$ws.Range("R12").Value = 30394.84920945689
$ws.Range("S12").Value = 0.1005744783936306
$ws.Range("T12").Value = 0.1005744783936306

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 122.909391
$ws.Range("H13").Value = 368.728173
$ws.Range("I13").Value = 0.4861499128584522
$ws.Range("J13").Value = 0.4861499128584522
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 13.164899
$ws.Range("N13").Value = 39.494697
$ws.Range("O13").Value = 0.09912031905407785
$ws.Range("P13").Value = 0.0991203190540779
$ws.Range("Q13").Value = 1618.089718666509
$ws.Range("R13").Value = 14562.80746799858
$ws.Range("S13").Value = 0.04818733447064192
$ws.Range("T13").Value = 0.04818733447064194

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.127183
$ws.Range("H14").Value = 45.381549
$ws.Range("I14").Value = 0.05983333443775553
$ws.Range("J14").Value = 0.05983333443775553
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 28.21795633333333
$ws.Range("N14").Value = 84.653869
$ws.Range("O14").Value = 0.2124568395711989
$ws.Range("P14").Value = 0.212456839571199
$ws.Range("Q14").Value = 426.8581893403423
$ws.Range("R14").Value = 3841.723704063081
$ws.Range("S14").Value = 0.01271200113565212
$ws.Range("T14").Value = 0.01271200113565212

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.127183
$ws.Range("H15").Value = 45.381549
$ws.Range("I15").Value = 0.05983333443775553
$ws.Range("J15").Value = 0.05983333443775553
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 63.95730733333333
$ws.Range("N15").Value = 191.871922
$ws.Range("O15").Value = 0.4815432848151522
$ws.Range("P15").Value = 0.4815432848151524
$ws.Range("Q15").Value = 967.4938922185753
$ws.Range("R15").Value = 8707.445029967177
$ws.Range("S15").Value = 0.02881234040660037
$ws.Range("T15").Value = 0.02881234040660037

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.127183
$ws.Range("H16").Value = 45.381549
$ws.Range("I16").Value = 0.05983333443775553
$ws.Range("J16").Value = 0.05983333443775553
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 27.47719633333334
$ws.Range("N16").Value = 82.431589
$ws.Range("O16").Value = 0.2068795565595709
$ws.Range("P16").Value = 0.2068795565595709
$ws.Range("Q16").Value = 415.6525772612624
$ws.Range("R16").Value = 3740.873195351361
$ws.Range("S16").Value = 0.01237829369596337
$ws.Range("T16").Value = 0.01237829369596337

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.127183
$ws.Range("H17").Value = 45.381549
$ws.Range("I17").Value = 0.05983333443775553
$ws.Range("J17").Value = 0.05983333443775553
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 13.164899
$ws.Range("N17").Value = 39.494697
$ws.Range("O17").Value = 0.09912031905407785
$ws.Range("P17").Value = 0.0991203190540779
$ws.Range("Q17").Value = 199.147836349517
$ws.Range("R17").Value = 1792.330527145653
$ws.Range("S17").Value = 0.005930699199539671
$ws.Range("T17").Value = 0.005930699199539674

Write-Host "Done updating cells"